# Remove three "Practice queries" slides from Session 3 (their content now
# lives in a separate notebook): the "Facebook Queries" slide (17), the
# "liberal students" practice-queries slide (18), and the "write down three
# queries" practice-queries slide (20). The REGEXP slide (19) and the
# "Go to the Facebook database" slide (21) are kept, sliding up to become
# the new slides 17 and 18.
#
# Delete from the highest index down so earlier deletions don't shift the
# indices of slides we still need to remove.
$p = $ppt.ActivePresentation
$p.Slides.Item(20).Delete()
$p.Slides.Item(18).Delete()
$p.Slides.Item(17).Delete()
